# The sheet's weekly data table (rows 68..153) needs 4 new rows of data
# inserted right before the current row 68, pushing the existing rows
# 68..153 down to 72..157 (Excel preserves their values/formatting
# automatically). We then fill in the values for the 4 brand-new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at 68, shifting rows 68:153 down to 72:157.
$ws.Rows("68:71").Insert()

# Row 68: Calameño / Extra
$ws.Range("A68").Value = 7
$ws.Range("B68").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C68").Value = 'Ñuble'
$ws.Range("D68").Value = 44571
$ws.Range("E68").Value = 16
$ws.Range("F68").Value = 100112027
$ws.Range("G68").Value = 'Melón'
$ws.Range("H68").Value = 'Calameño'
$ws.Range("I68").Value = 'Extra'
$ws.Range("J68").Value = 1500
$ws.Range("K68").Value = 700
$ws.Range("L68").Value = 700
$ws.Range("M68").Value = 700
$ws.Range("N68").Value = '$/unidad'
$ws.Range("O68").Value = 'Región del Maule'
$ws.Range("P68").Value = 700
$ws.Range("Q68").Value = 1
$ws.Range("R68").Value = 'Hortaliza'

# Row 69: Calameño / Primera
$ws.Range("A69").Value = 7
$ws.Range("B69").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C69").Value = 'Ñuble'
$ws.Range("D69").Value = 44571
$ws.Range("E69").Value = 16
$ws.Range("F69").Value = 100112027
$ws.Range("G69").Value = 'Melón'
$ws.Range("H69").Value = 'Calameño'
$ws.Range("I69").Value = 'Primera'
$ws.Range("J69").Value = 2400
$ws.Range("K69").Value = 500
$ws.Range("L69").Value = 600
$ws.Range("M69").Value = 550
$ws.Range("N69").Value = '$/unidad'
$ws.Range("O69").Value = 'Región del Maule'
$ws.Range("P69").Value = 550
$ws.Range("Q69").Value = 1
$ws.Range("R69").Value = 'Hortaliza'

# Row 70: Tuna / Extra
$ws.Range("A70").Value = 7
$ws.Range("B70").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C70").Value = 'Ñuble'
$ws.Range("D70").Value = 44571
$ws.Range("E70").Value = 16
$ws.Range("F70").Value = 100112027
$ws.Range("G70").Value = 'Melón'
$ws.Range("H70").Value = 'Tuna'
$ws.Range("I70").Value = 'Extra'
$ws.Range("J70").Value = 1500
$ws.Range("K70").Value = 700
$ws.Range("L70").Value = 700
$ws.Range("M70").Value = 700
$ws.Range("N70").Value = '$/unidad'
$ws.Range("O70").Value = 'Región del Maule'
$ws.Range("P70").Value = 700
$ws.Range("Q70").Value = 1
$ws.Range("R70").Value = 'Hortaliza'

# Row 71: Tuna / Primera
$ws.Range("A71").Value = 7
$ws.Range("B71").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C71").Value = 'Ñuble'
$ws.Range("D71").Value = 44571
$ws.Range("E71").Value = 16
$ws.Range("F71").Value = 100112027
$ws.Range("G71").Value = 'Melón'
$ws.Range("H71").Value = 'Tuna'
$ws.Range("I71").Value = 'Primera'
$ws.Range("J71").Value = 2400
$ws.Range("K71").Value = 500
$ws.Range("L71").Value = 600
$ws.Range("M71").Value = 550
$ws.Range("N71").Value = '$/unidad'
$ws.Range("O71").Value = 'Región del Maule'
$ws.Range("P71").Value = 550
$ws.Range("Q71").Value = 1
$ws.Range("R71").Value = 'Hortaliza'
